$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 to 6 (NOUBAIL MOUNTASSIR, NOUBAIL MOHAMMED, totals row)
$ws.Range("A4:K6").EntireRow.Delete()

# Update row 2 with new data
$ws.Range("A2").Value = "NASIRI HASNAA"
$ws.Range("B2").Value = ""
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "546576878798989898090090"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "CIH"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "901/CASABLANCA"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 8000.01
$ws.Range("J2").Value = 800.01
$ws.Range("K2").Value = 7200

# Update row 3 with blank text cells but numeric totals
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = " "
$ws.Range("I3").Value = 8000.01
$ws.Range("J3").Value = 800.01
$ws.Range("K3").Value = 7200
